# Apply the PGP wind/solar case_input edit:
# Insert a new row at row 46 (a NORMALIZE_DEMAND_TO_ONE flag row), pushing all
# rows from the old row 46 onward down by one, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing row 46 (copies formatting/styles along with it) and
# insert the copy above it, shifting rows 46:106 down to 47:107.
$ws.Rows("46:46").Copy()
$ws.Rows("46:46").Insert()

# Overwrite the new row 46 with the NORMALIZE_DEMAND_TO_ONE flag values.
$ws.Range("A46").Value = "NORMALIZE_DEMAND_TO_ONE"
$ws.Range("B46").Value = $true
$ws.Range("C46").Value = "Nomalize demand 1."

# Update the view: move the active selection to C46.
$ws.Range("C46").Select()
